# Update last_edited_time stamps and a few numeric properties on
# worksheet LUY_KE_THANG_SOC_TRANG to account for differences between
# mac and win runs (per commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LUY_KE_THANG_SOC_TRANG")

$oldTime = "2024-07-19T12:51:00.000Z"
$newTime = "2024-07-20T13:34:00.000Z"

# Rows whose "last_edited_time" (column D) needs to be bumped.
$rowsToStamp = @(2, 3, 6, 8, 11, 13)
foreach ($r in $rowsToStamp) {
    $cell = $ws.Range("D$r")
    if ($cell.Value2 -eq $oldTime) {
        $cell.Value = $newTime
    }
}

# Numeric property updates on row 13 (Thang 7).
$ws.Range("S13").Value = 143722000
$ws.Range("W13").Value = 19866000
$ws.Range("AE13").Value = 163588000
$ws.Range("AN13").Value = 26000000
